$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All of these columns (D = Price, E = Volume(1h))
# store plain text that looks numeric, so force text format ("@") before assignment
# to avoid Excel auto-converting the strings into numeric/percentage values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.81%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.12%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.167"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.38%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07497"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.86%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.789"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.43%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.795"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.38%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.667"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.49%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9257"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.57%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1713"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.78%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07545"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.93%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07947"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.36%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.79%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09889"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.27%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001496"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.93%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04663"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.31%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006507"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "5.46%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.462"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.50%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.44%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.60%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.15%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.563"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.47%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1549"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.27%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.88%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004415"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.70%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "19.81%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001808"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "8.72%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01658"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.72%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04546"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.08%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006996"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-4.13%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.74%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002059"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.73%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01313"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.05%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006089"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.01%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.930"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.98%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-5.62%"
